# Update "想去人数" (people interested) counts across sheets, as published
# by the gh-pages data refresh at commit 456a3b4.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 37635
$ws.Range("F4").Value = 636
$ws.Range("F5").Value = 768
$ws.Range("F6").Value = 477
$ws.Range("F7").Value = 364
$ws.Range("F10").Value = 94
$ws.Range("F11").Value = 712
$ws.Range("F12").Value = 542
$ws.Range("F13").Value = 43
$ws.Range("F15").Value = 18
$ws.Range("F16").Value = 648
$ws.Range("F17").Value = 177
$ws.Range("F19").Value = 442
$ws.Range("F20").Value = 1168
$ws.Range("F22").Value = 828
$ws.Range("F23").Value = 2528
$ws.Range("F24").Value = 1007
$ws.Range("F26").Value = 108
$ws.Range("F27").Value = 1162
$ws.Range("F29").Value = 774
$ws.Range("F31").Value = 1158

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 400

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 630

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 630
$ws.Range("F3").Value = 37635
$ws.Range("F5").Value = 636
$ws.Range("F6").Value = 768
$ws.Range("F7").Value = 477
$ws.Range("F9").Value = 364
$ws.Range("F11").Value = 400
$ws.Range("F16").Value = 94
$ws.Range("F17").Value = 712
$ws.Range("F18").Value = 542
$ws.Range("F19").Value = 43
$ws.Range("F25").Value = 18
$ws.Range("F27").Value = 648
$ws.Range("F28").Value = 177
$ws.Range("F30").Value = 442
$ws.Range("F31").Value = 1168
$ws.Range("F33").Value = 828
$ws.Range("F34").Value = 2528
$ws.Range("F35").Value = 1007
$ws.Range("F37").Value = 108
$ws.Range("F38").Value = 1162
$ws.Range("F41").Value = 774
$ws.Range("F43").Value = 1158
